$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2021" column (O) mirroring column N's layout/styles/formatting.

# Copy formatting from column N (rows 3-10) into column O first.
$ws.Range("N3:N10").Copy()
$ws.Range("O3:O10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 4 header: year 2021
$ws.Range("O4").Value = 2021

# Data rows
$ws.Range("O6").Value = 1860
$ws.Range("O7").Value = 1
$ws.Range("O8").Value = 510
$ws.Range("O9").Value = 178
$ws.Range("O10").Value = 821

# Update selection to match the post-edit state
$ws.Range("P9").Select()
